$d = $word.ActiveDocument

# 1. First paragraph: add "Assunto: Remembramento" before the tabs, and "3" before "ª análise"
$d.Paragraphs.Item(1).Range.Text = "Assunto: Remembramento				3ª análise"

# 2. Solicitação de demanda:
$d.Paragraphs.Item(2).Range.Text = "Solicitação de demanda: Telefone"

# 3. Contribuinte:
$d.Paragraphs.Item(3).Range.Text = "Contribuinte: Carlos"

# 4. Inscrição Imobiliária:
$d.Paragraphs.Item(4).Range.Text = "Inscrição Imobiliária: 123456789012345"

# 5. Endereço do imóvel:
$d.Paragraphs.Item(5).Range.Text = "Endereço do imóvel: Rua Cibele, nº 145 - bairro Vila Amélia, Itaira - MG"

# 6. Dados recebidos: (adds literal newline characters inside the same run/text)
$d.Paragraphs.Item(6).Range.Text = "Dados recebidos: - Compra e venda`n 1234`n"

# 7. Final empty paragraph gets the conclusion text
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = "Após verificação dos arquivos apresentados à Prefeitura Municipal de Itabira referentes ao levantamento realizado, não foram identificados deslocamentos, sobreposições, nem invasão de vias públicas. Recomenda-se que a Prefeitura Municipal de Itabira opte pelo deferimento do processo XXXX/XX/XXXX."
